$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.911.02"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3
$ws.Range("D3").Value = "1.894.97"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'0.7736"
$ws.Range("E5").Value = "  -2.43%  "

# Row 6
$ws.Range("D6").Value = "'244.69"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.3137"
$ws.Range("E8").Value = "  -0.79%  "

# Row 9
$ws.Range("D9").Value = "'25.70"
$ws.Range("E9").Value = "  +1.15%  "

# Row 10
$ws.Range("E10").Value = "  +0.53%  "

# Row 11
$ws.Range("D11").Value = "'0.08909"
$ws.Range("E11").Value = "  +9.91%  "

# Row 12
$ws.Range("D12").Value = "'0.7729"
$ws.Range("E12").Value = "  +0.74%  "

# Row 13
$ws.Range("D13").Value = "'5.445"
$ws.Range("E13").Value = "  -2.67%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.844.43"
$ws.Range("E14").Value = "  -1.28%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'94.46"
$ws.Range("E15").Value = "  +1.95%  "

# Row 16
$ws.Range("D16").Value = "'6.195"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").Value = "29.863.87"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").Value = "'13.98"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19
$ws.Range("D19").Value = "'246.13"
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").Value = "'0.000007876"
$ws.Range("E20").Value = "  +0.96%  "

# Row 21
$ws.Range("D21").Value = "'8.133"
$ws.Range("E21").Value = "  -1.00%  "

# Row 22
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").Value = "2.115.04"
$ws.Range("E23").Value = "  -1.25%  "

# Row 25
$ws.Range("D25").Value = "'0.1598"
$ws.Range("E25").Value = "  -4.52%  "

# Row 26
$ws.Range("D26").Value = "'9.532"
$ws.Range("E26").Value = "  +0.91%  "

# Row 27
$ws.Range("D27").Value = "'162.96"
$ws.Range("E27").Value = "  -0.80%  "

# Row 28
$ws.Range("D28").Value = "'18.82"
$ws.Range("E28").Value = "  +0.51%  "

# Row 29
$ws.Range("D29").Value = "'2.047"
$ws.Range("E29").Value = "  -0.89%  "

# Row 30
$ws.Range("E30").Value = "  +1.82%  "

# Row 31
$ws.Range("E31").Value = "  -0.29%  "

# Row 32
$ws.Range("D32").Value = "'4.543"
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("D33").Value = "'4.108"
$ws.Range("E33").Value = "  +0.32%  "

# Row 34
$ws.Range("D34").Value = "'0.05520"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35
$ws.Range("D35").Value = "'1.248"
$ws.Range("E35").Value = "  -2.60%  "

# Row 36
$ws.Range("D36").Value = "'0.7519"
$ws.Range("E36").Value = "  +1.51%  "

# Row 37
$ws.Range("D37").Value = "'0.9959"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38
$ws.Range("D38").Value = "'2.718"
$ws.Range("E38").Value = "  +3.39%  "

# Row 39
$ws.Range("D39").Value = "'0.01962"
$ws.Range("E39").Value = "  +1.63%  "

# Row 40
$ws.Range("D40").Value = "'2.792"
$ws.Range("E40").Value = "  +0.36%  "

# Row 41
$ws.Range("D41").Value = "'0.4506"
$ws.Range("E41").Value = "  +1.96%  "

# Row 42
$ws.Range("D42").Value = "'73.91"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43
$ws.Range("D43").Value = "'6.050"
$ws.Range("E43").Value = "  +2.53%  "

# Row 44
$ws.Range("D44").Value = "1.086.15"
$ws.Range("E44").Value = "  -6.14%  "

# Row 45
$ws.Range("D45").Value = "'0.8550"
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "'1.890"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48
$ws.Range("D48").Value = "'102.61"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.865"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.608"
$ws.Range("E50").Value = "  +1.99%  "

# Row 51
$ws.Range("D51").Value = "'2.993"
$ws.Range("E51").Value = "  -1.58%  "
